$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: new estimate for InfAV, using ATV between y-1 and y (instead of m-1/m or q-1/q)
$ws.Range("B2").Value = [double]"-3.973643103449831E-08"

# Row 3: newly populated moments for this column
$ws.Range("B3").Value = 0
$ws.Range("C3").Value = 2.578947368421053
$ws.Range("D3").Value = 2.321052631578948
$ws.Range("E3").Value = -0.03515381148675916
$ws.Range("F3").Value = 0.6375222362110279
$ws.Range("G3").Value = 0.3949431464863624
$ws.Range("H3").Value = 0.3886891279823181
$ws.Range("I3").Value = 0.1047529638783499
$ws.Range("J3").Value = 0.08627647523793786
$ws.Range("K3").Value = 1.360294011187797
$ws.Range("L3").Value = 0.0007933425278022433
$ws.Range("M3").Value = 0.0005139105411266024

# Row 4: newly populated moments for this column
$ws.Range("B4").Value = 0
$ws.Range("C4").Value = 0.252525252525252
$ws.Range("D4").Value = 0.247474747474747
$ws.Range("E4").Value = -0.006323887507055811
$ws.Range("F4").Value = 0.01301996818610838
$ws.Range("G4").Value = 0.008470320293531447
$ws.Range("H4").Value = 0.01756716503656709
$ws.Range("I4").Value = 0.0001789228060869995
$ws.Range("J4").Value = 0.000164526812178966
$ws.Range("K4").Value = 0.04302632986889662
$ws.Range("L4").Value = [double]"1.709752113614707E-06"
$ws.Range("M4").Value = [double]"1.247776857399013E-06"

# Row 5: newly populated moments for this column
$ws.Range("B5").Value = 0
$ws.Range("C5").Value = 3.120526315789474
$ws.Range("D5").Value = 2.808473684210527
$ws.Range("E5").Value = -0.02608882149306242
$ws.Range("F5").Value = 0.4172501193485463
$ws.Range("G5").Value = 0.1265177469816006
$ws.Range("H5").Value = 0.3498577899066002
$ws.Range("I5").Value = 0.09382889719047643
$ws.Range("J5").Value = 0.05847461944715957
$ws.Range("K5").Value = 1.170398281091975
$ws.Range("L5").Value = 0.0006330841605619485
$ws.Range("M5").Value = 0.0002861148571598649

# Row 6: newly populated moments for this column
$ws.Range("B6").Value = 0
$ws.Range("C6").Value = 2.578947368421053
$ws.Range("D6").Value = 2.321052631578948
$ws.Range("E6").Value = 0.005770085637884103
$ws.Range("F6").Value = 0.5122349381704008
$ws.Range("G6").Value = 0.01618863295228209
$ws.Range("H6").Value = [double]"9.987853222029977E-29"
$ws.Range("I6").Value = [double]"4.781103680424108E-56"
$ws.Range("J6").Value = [double]"2.017022683526623E-56"
$ws.Range("K6").Value = 0.4900000000000072
$ws.Range("L6").Value = [double]"1.597443333072549E-29"
$ws.Range("M6").Value = [double]"1.600871323057683E-29"

# Row 7: newly populated moments for this column
$ws.Range("B7").Value = 0
$ws.Range("C7").Value = 2.578947368421053
$ws.Range("D7").Value = 2.321052631578948
$ws.Range("E7").Value = 0.005770085637884103
$ws.Range("F7").Value = 0.5122349381704008
$ws.Range("G7").Value = 0.01618863295228209
$ws.Range("H7").Value = [double]"9.987853222029977E-29"
$ws.Range("I7").Value = [double]"4.781103680424108E-56"
$ws.Range("J7").Value = [double]"2.017022683526623E-56"
$ws.Range("K7").Value = 0.4900000000000072
$ws.Range("L7").Value = [double]"1.597443333072549E-29"
$ws.Range("M7").Value = [double]"1.600871323057683E-29"

# Row 8: newly populated moments for this column
$ws.Range("B8").Value = 0
$ws.Range("C8").Value = 3.120526315789474
$ws.Range("D8").Value = 2.808473684210527
$ws.Range("E8").Value = 0.1135701804545673
$ws.Range("F8").Value = 0.8984155410771357
$ws.Range("G8").Value = 0.5818480429174054
$ws.Range("H8").Value = 0.210107003145396
$ws.Range("I8").Value = 0.0002076668520766055
$ws.Range("J8").Value = 0.000105206809290575
$ws.Range("K8").Value = 1.507615341825497
$ws.Range("L8").Value = [double]"5.048709793414476E-29"
$ws.Range("M8").Value = [double]"5.05954393460206E-29"

# Row 9: newly populated moments for this column
$ws.Range("B9").Value = 0
$ws.Range("C9").Value = 2.578947368421053
$ws.Range("D9").Value = 2.321052631578948
$ws.Range("E9").Value = 0.005646565565627371
$ws.Range("F9").Value = 1.574635842517005
$ws.Range("G9").Value = -0.1394394934465455
$ws.Range("H9").Value = 1.626654514699598
$ws.Range("I9").Value = 4.73008038487975
$ws.Range("J9").Value = -0.342206364509281
$ws.Range("K9").Value = 0.4900000000000072
$ws.Range("L9").Value = [double]"1.597443333072549E-29"
$ws.Range("M9").Value = [double]"1.600871323057683E-29"

# Row 10: newly populated moments for this column
$ws.Range("B10").Value = 0
$ws.Range("C10").Value = 2.578947368421053
$ws.Range("D10").Value = 2.321052631578948
$ws.Range("E10").Value = 0.005648640377334755
$ws.Range("F10").Value = 1.574872049104573
$ws.Range("G10").Value = -0.1394606933301283
$ws.Range("H10").Value = 0.372711814749233
$ws.Range("I10").Value = 0.2483271448891567
$ws.Range("J10").Value = -0.01796568399411001
$ws.Range("K10").Value = 0.4900000000000072
$ws.Range("L10").Value = [double]"1.597443333072549E-29"
$ws.Range("M10").Value = [double]"1.600871323057683E-29"

# Row 11: newly populated moments for this column
$ws.Range("B11").Value = 0
$ws.Range("C11").Value = 2.578947368421053
$ws.Range("D11").Value = 2.321052631578948
$ws.Range("E11").Value = 0.004364864289895591
$ws.Range("F11").Value = 1.432229617436627
$ws.Range("G11").Value = -0.1265367874246578
$ws.Range("H11").Value = 0.3367230746719863
$ws.Range("I11").Value = 0.2026859441625896
$ws.Range("J11").Value = -0.01466368738906204
$ws.Range("K11").Value = 0.4900000000000072
$ws.Range("L11").Value = [double]"1.597443333072549E-29"
$ws.Range("M11").Value = [double]"1.600871323057683E-29"

# Row 12: newly populated moments for this column
$ws.Range("B12").Value = 0
$ws.Range("C12").Value = 2.292631578947369
$ws.Range("D12").Value = 2.063368421052632
$ws.Range("E12").Value = 0.01289490163026139
$ws.Range("F12").Value = 2.195066048322307
$ws.Range("G12").Value = -1.250146626724897
$ws.Range("H12").Value = [double]"1.433553253389728E-07"
$ws.Range("I12").Value = [double]"7.220022858142143E-17"
$ws.Range("J12").Value = [double]"4.48647694133143E-17"
$ws.Range("K12").Value = 0.4359236885998423
$ws.Range("L12").Value = [double]"2.415886522239349E-30"
$ws.Range("M12").Value = [double]"2.421070828081064E-30"

# Row 13: newly populated moments for this column
$ws.Range("B13").Value = 0
$ws.Range("C13").Value = 2.321243523316065
$ws.Range("D13").Value = 2.15875647668394
$ws.Range("E13").Value = 0.02247094754209927
$ws.Range("F13").Value = 0.6574192142455283
$ws.Range("G13").Value = -0.08282831512109545
$ws.Range("H13").Value = 0.04939396741652854
$ws.Range("I13").Value = [double]"9.027236561935378E-06"
$ws.Range("J13").Value = [double]"1.98184864247584E-06"
$ws.Range("K13").Value = 0.3317145783344422
$ws.Range("L13").Value = [double]"3.993608332681372E-30"
$ws.Range("M13").Value = [double]"4.002178307644208E-30"

# Row 14: newly populated moments for this column
$ws.Range("B14").Value = 0
$ws.Range("C14").Value = 0.07038976148923795
$ws.Range("D14").Value = 0.06405468295520654
$ws.Range("E14").Value = 0.03179736183752695
$ws.Range("F14").Value = 0.01448287415193935
$ws.Range("G14").Value = 0.006944976028680937
$ws.Range("H14").Value = 0.01138370561386528
$ws.Range("I14").Value = [double]"5.952666499984613E-07"
$ws.Range("J14").Value = [double]"2.491922779909072E-07"
$ws.Range("K14").Value = 0.03432440880166791
$ws.Range("L14").Value = [double]"3.510007323645737E-32"
$ws.Range("M14").Value = [double]"3.517539528202917E-32"
